$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells in row 1, copying the style used by the other headers (e.g. A1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the win/loss/tie record for every data row (2 through 46)
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 99   # column AD
    $ws.Cells.Item($r, 31).Value = 62   # column AE
    $ws.Cells.Item($r, 32).Value = 0    # column AF
}
